$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (AstronautDistressCall) status moved from "Voice Recording" to
# "Reaper Done" - copy F4's formatting (style "40% - Accent2", used by the
# other "Reaper Done" status cell) onto F11, then update its text.
$ws.Range("F4").Copy()
$ws.Range("F11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F11").Value2 = "Reaper Done"

# Legend / KEY section: rename the "Voice Recording" entry to
# "(DONE) Voice Recording", keeping the bold lead-in run and the existing
# plain " - script/voice work needed" description run intact.
$legend = $ws.Range("D18")
$legend.Characters(1, 15).Text = "(DONE) Voice Recording"
$legend.Characters(1, 22).Font.Bold = $true
$total = $legend.Characters().Count
$legend.Characters(23, $total - 22).Font.Name = "Calibri"
$legend.Characters(23, $total - 22).Font.Size = 11

# Restore the active cell/selection to F10 (last selected cell in the edit).
[void]$ws.Range("F10").Select()

$excel.CutCopyMode = $false
